$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 11980
$ws.Range("I19").Value = 9900
$ws.Range("J19").Value = 12500
$ws.Range("K19").Value = 9900
$ws.Range("L19").Value = 12500
$ws.Range("M19").Value = -9725
$ws.Range("N19").Value = -12850
# Row 43
$ws.Range("H43").Value = 6941.8887
$ws.Range("I43").Value = 2608.3333
$ws.Range("J43").Value = 9108.666999999999
$ws.Range("K43").Value = 2608.3333
$ws.Range("L43").Value = 9108.666999999999
$ws.Range("M43").Value = -2539.3333
$ws.Range("N43").Value = -9246.666999999999
# Row 76
$ws.Range("H76").Value = 3153.195
$ws.Range("I76").Value = 3119.4517
$ws.Range("J76").Value = 3257.8
$ws.Range("K76").Value = 3119.4517
$ws.Range("L76").Value = 3257.8
$ws.Range("M76").Value = -2804.4517
$ws.Range("N76").Value = -3887.8
# Row 79
$ws.Range("H79").Value = 3153.195
$ws.Range("I79").Value = 3119.4517
$ws.Range("J79").Value = 3257.8
$ws.Range("K79").Value = 3119.4517
$ws.Range("L79").Value = 3257.8
$ws.Range("M79").Value = -2027.4517
$ws.Range("N79").Value = -5441.8
# Row 107
$ws.Range("H107").Value = 5231.75
$ws.Range("I107").Value = 9279.454
$ws.Range("J107").Value = 284.55554
$ws.Range("K107").Value = 9279.454
$ws.Range("L107").Value = 284.55554
$ws.Range("M107").Value = -7359.454
$ws.Range("N107").Value = -4124.55554
# Row 137
$ws.Range("H137").Value = 40100600
$ws.Range("I137").Value = 999.6667
$ws.Range("K137").Value = 2999.0001
$ws.Range("M137").Value = -449.0001000000002
# Row 141
$ws.Range("H141").Value = 2831.25
$ws.Range("I141").Value = 2124.1667
$ws.Range("J141").Value = 4952.5
$ws.Range("K141").Value = 6372.500100000001
$ws.Range("L141").Value = 14857.5
$ws.Range("M141").Value = -1192.500100000001
$ws.Range("N141").Value = -25217.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3540.337
$ws.Range("I32").Value = 2914.1208
$ws.Range("J32").Value = 17786.75
$ws.Range("K32").Value = 2914.1208
$ws.Range("L32").Value = 17786.75
$ws.Range("M32").Value = -2627.1208
$ws.Range("N32").Value = -18360.75
# Row 63
$ws.Range("H63").Value = 16130858
$ws.Range("I63").Value = 16668440
$ws.Range("K63").Value = 16668440
$ws.Range("M63").Value = -16667754
# Row 66
$ws.Range("H66").Value = 16130858
$ws.Range("I66").Value = 16668440
$ws.Range("K66").Value = 83342200
$ws.Range("M66").Value = -83338768
# Row 88
$ws.Range("H88").Value = 1966.6666
$ws.Range("I88").Value = 1500
$ws.Range("J88").Value = 2200
$ws.Range("K88").Value = 1500
$ws.Range("L88").Value = 2200
$ws.Range("M88").Value = -1094
$ws.Range("N88").Value = -3012
# Row 91
$ws.Range("H91").Value = 1966.6666
$ws.Range("I91").Value = 1500
$ws.Range("J91").Value = 2200
$ws.Range("K91").Value = 1500
$ws.Range("L91").Value = 2200
$ws.Range("M91").Value = -96
$ws.Range("N91").Value = -5008
# Row 106
$ws.Range("H106").Value = 50370
$ws.Range("J106").Value = 50370
$ws.Range("L106").Value = 50370
$ws.Range("N106").Value = -52894

$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 25707.125
$ws.Range("I80").Value = 40370
$ws.Range("K80").Value = 40370
$ws.Range("M80").Value = -39372
# Row 83
$ws.Range("H83").Value = 25707.125
$ws.Range("I83").Value = 40370
$ws.Range("K83").Value = 201850
$ws.Range("M83").Value = -196858
# Row 86
$ws.Range("H86").Value = 20835874
$ws.Range("I86").Value = 2353.1875
$ws.Range("J86").Value = 62502910
$ws.Range("K86").Value = 2353.1875
$ws.Range("L86").Value = 62502910
$ws.Range("M86").Value = -1230.1875
$ws.Range("N86").Value = -62505156
# Row 89
$ws.Range("H89").Value = 20835874
$ws.Range("I89").Value = 2353.1875
$ws.Range("J89").Value = 62502910
$ws.Range("K89").Value = 11765.9375
$ws.Range("L89").Value = 312514550
$ws.Range("M89").Value = -6149.9375
$ws.Range("N89").Value = -312525782

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5850235.5
$ws.Range("I31").Value = 1799.7297
$ws.Range("J31").Value = 16669841
$ws.Range("K31").Value = 1799.7297
$ws.Range("L31").Value = 16669841
$ws.Range("M31").Value = -1504.7297
$ws.Range("N31").Value = -16670431
# Row 34
$ws.Range("H34").Value = 5850235.5
$ws.Range("I34").Value = 1799.7297
$ws.Range("J34").Value = 16669841
$ws.Range("K34").Value = 1799.7297
$ws.Range("L34").Value = 16669841
$ws.Range("M34").Value = -1597.7297
$ws.Range("N34").Value = -16670245
# Row 70
$ws.Range("H70").Value = 25178.889
$ws.Range("J70").Value = 25178.889
$ws.Range("L70").Value = 25178.889
$ws.Range("N70").Value = -25808.889
# Row 73
$ws.Range("H73").Value = 25178.889
$ws.Range("J73").Value = 25178.889
$ws.Range("L73").Value = 25178.889
$ws.Range("N73").Value = -27362.889
# Row 75
$ws.Range("H75").Value = 42196
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 42196
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 42196
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -44192
# Row 78
$ws.Range("H78").Value = 42196
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 42196
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 126588
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -136572
# Row 80
$ws.Range("H80").Value = 24000
$ws.Range("J80").Value = 24000
$ws.Range("L80").Value = 24000
$ws.Range("N80").Value = -26246
# Row 83
$ws.Range("H83").Value = 24000
$ws.Range("J83").Value = 24000
$ws.Range("L83").Value = 72000
$ws.Range("N83").Value = -83232
# Row 86
$ws.Range("H86").Value = 2713.9678
$ws.Range("I86").Value = 2528.0527
$ws.Range("J86").Value = 3008.3333
$ws.Range("K86").Value = 2528.0527
$ws.Range("L86").Value = 3008.3333
$ws.Range("M86").Value = -1405.0527
$ws.Range("N86").Value = -5254.3333
# Row 89
$ws.Range("H89").Value = 2713.9678
$ws.Range("I89").Value = 2528.0527
$ws.Range("J89").Value = 3008.3333
$ws.Range("K89").Value = 12640.2635
$ws.Range("L89").Value = 15041.6665
$ws.Range("M89").Value = -7024.263500000001
$ws.Range("N89").Value = -26273.6665
# Row 134
$ws.Range("H134").Value = 918254.4
$ws.Range("I134").Value = 2952.8572
$ws.Range("J134").Value = 1986106.1
$ws.Range("K134").Value = 8858.571599999999
$ws.Range("L134").Value = 5958318.300000001
$ws.Range("M134").Value = -6323.571599999999
$ws.Range("N134").Value = -5963388.300000001

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 8010.479
$ws.Range("I70").Value = 9325.723
$ws.Range("J70").Value = 4064.75
$ws.Range("K70").Value = 9325.723
$ws.Range("L70").Value = 4064.75
$ws.Range("M70").Value = -9055.723
$ws.Range("N70").Value = -4604.75
# Row 73
$ws.Range("H73").Value = 8010.479
$ws.Range("I73").Value = 9325.723
$ws.Range("J73").Value = 4064.75
$ws.Range("K73").Value = 9325.723
$ws.Range("L73").Value = 4064.75
$ws.Range("M73").Value = -8389.723
$ws.Range("N73").Value = -5936.75
# Row 86
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
# Row 89
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
# Row 132
$ws.Range("H132").Value = 4387.0225
$ws.Range("I132").Value = 5206.871
$ws.Range("K132").Value = 15620.613
$ws.Range("M132").Value = -13090.613

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5304.4136
$ws.Range("I7").Value = 5253.4287
$ws.Range("J7").Value = 5352
$ws.Range("K7").Value = 5253.4287
$ws.Range("L7").Value = 5352
$ws.Range("M7").Value = -5141.4287
$ws.Range("N7").Value = -5576
# Row 126
$ws.Range("H126").Value = 5304.4136
$ws.Range("I126").Value = 5253.4287
$ws.Range("J126").Value = 5352
$ws.Range("K126").Value = 15760.2861
$ws.Range("L126").Value = 16056
$ws.Range("M126").Value = -13290.2861
$ws.Range("N126").Value = -20996

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1976.1471
$ws.Range("I122").Value = 2004.7273
$ws.Range("J122").Value = 1923.75
$ws.Range("K122").Value = 6014.1819
$ws.Range("L122").Value = 5771.25
$ws.Range("M122").Value = -3564.1819
$ws.Range("N122").Value = -10671.25
# Row 126
$ws.Range("H126").Value = 1593.4667
$ws.Range("I126").Value = 1242
$ws.Range("J126").Value = 2999.3333
$ws.Range("K126").Value = 3726
$ws.Range("L126").Value = 8997.999899999999
$ws.Range("M126").Value = -1256
$ws.Range("N126").Value = -13937.9999
# Row 136
$ws.Range("H136").Value = 1095.5714
$ws.Range("I136").Value = 1208.8636
$ws.Range("K136").Value = 3626.5908
$ws.Range("M136").Value = -1076.5908
